$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.949.84"
$ws.Range("E2").Value = "  -3.45%  "
$ws.Range("D3").Value = "2.921.58"
$ws.Range("E3").Value = "  -3.95%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("E5").Value = "  -1.52%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "145.50"
$ws.Range("E6").Value = "  -5.87%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.505"
$ws.Range("E8").Value = "  -2.31%  "
$ws.Range("D9").Value = "2.920.46"
$ws.Range("E9").Value = "  -3.96%  "
$ws.Range("E10").Value = "  -1.16%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.144"
$ws.Range("E11").Value = "  -5.31%  "
$ws.Range("E12").Value = "  -3.58%  "
$ws.Range("E13").Value = "  -3.90%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "33.60"
$ws.Range("E14").Value = "  -6.86%  "
$ws.Range("E15").Value = "  +0.34%  "
$ws.Range("D16").Value = "3.404.63"
$ws.Range("E16").Value = "  -3.98%  "
$ws.Range("D17").Value = "60.893.51"
$ws.Range("E17").Value = "  -3.47%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.76"
$ws.Range("E18").Value = "  -5.04%  "
$ws.Range("D19").Value = "2.921.42"
$ws.Range("E19").Value = "  -4.00%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "428.56"
$ws.Range("E20").Value = "  -5.97%  "
$ws.Range("E21").Value = "  -4.95%  "
$ws.Range("E22").Value = "  -2.49%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.12"
$ws.Range("E23").Value = "  -5.70%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "80.75"
$ws.Range("E24").Value = "  -2.87%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.24"
$ws.Range("E25").Value = "  -3.10%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "10.76"
$ws.Range("E26").Value = "  -4.91%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.01"
$ws.Range("E27").Value = "  -3.01%  "
$ws.Range("E28").Value = "  +0.01%  "
$ws.Range("E29").Value = "  +0.01%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.22"
$ws.Range("E30").Value = "  -3.18%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.61"
$ws.Range("E31").Value = "  -3.28%  "
$ws.Range("E32").Value = "  -3.31%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "26.66"
$ws.Range("E33").Value = "  -3.91%  "
$ws.Range("E34").Value = "  -3.39%  "
$ws.Range("D35").Value = "0.0₃0870"
$ws.Range("E35").Value = "  +0.90%  "
$ws.Range("E36").Value = "  -2.64%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.65"
$ws.Range("E37").Value = "  -5.12%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.01"
$ws.Range("E38").Value = "  -5.62%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.127"
$ws.Range("E39").Value = "  -3.79%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "49.63"
$ws.Range("E40").Value = "  -1.64%  "
$ws.Range("E41").Value = "  -5.34%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.65"
$ws.Range("E42").Value = "  -5.39%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.298"
$ws.Range("E43").Value = "  -2.17%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "40.96"
$ws.Range("E44").Value = "  -6.05%  "
$ws.Range("B45").Value = "Bittensor"
$ws.Range("C45").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "379.92"
$ws.Range("E45").Value = "  -3.52%  "
$ws.Range("B46").Value = "VeChain"
$ws.Range("C46").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0352"
$ws.Range("E46").Value = "  -2.80%  "
$ws.Range("D47").Value = "2.689.98"
$ws.Range("E47").Value = "  -1.38%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "132.60"
$ws.Range("E48").Value = "  +0.14%  "
$ws.Range("E49").Value = "  -0.05%  "
$ws.Range("E50").Value = "  -0.35%  "
$ws.Range("E51").Value = "  -2.33%  "
